$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$cs = $nm.ColorScheme
Write-Output (Get-Member -InputObject $cs)
